$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 7719.2954
$ws.Range("I132").Value = 7118.3
$ws.Range("J132").Value = 9007.143
$ws.Range("K132").Value = 21354.9
$ws.Range("L132").Value = 27021.429
$ws.Range("M132").Value = -18824.9
$ws.Range("N132").Value = -32081.429

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1940.4773
$ws.Range("I61").Value = 1909.7241
$ws.Range("J61").Value = 1999.9333
$ws.Range("K61").Value = 1909.7241
$ws.Range("L61").Value = 1999.9333
$ws.Range("M61").Value = -1697.7241
$ws.Range("N61").Value = -2423.9333

# Row 122
$ws.Range("H122").Value = 1388.75
$ws.Range("I122").Value = 1230
$ws.Range("J122").Value = 1484
$ws.Range("K122").Value = 3690
$ws.Range("L122").Value = 4452
$ws.Range("M122").Value = -1240
$ws.Range("N122").Value = -9352

# Row 136
$ws.Range("H136").Value = 1940.4773
$ws.Range("I136").Value = 1909.7241
$ws.Range("J136").Value = 1999.9333
$ws.Range("K136").Value = 5729.1723
$ws.Range("L136").Value = 5999.7999
$ws.Range("M136").Value = -3179.1723
$ws.Range("N136").Value = -11099.7999

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1587.6364
$ws.Range("I20").Value = 879.3
$ws.Range("J20").Value = 2677.3845
$ws.Range("K20").Value = 879.3
$ws.Range("L20").Value = 2677.3845
$ws.Range("M20").Value = -632.3
$ws.Range("N20").Value = -3171.3845

# Row 123
$ws.Range("H123").Value = 48000
$ws.Range("J123").Value = 48000
$ws.Range("L123").Value = 48000
$ws.Range("N123").Value = -57800

# Row 134
$ws.Range("H134").Value = 5077.477
$ws.Range("I134").Value = 2329.5386
$ws.Range("J134").Value = 6229.839
$ws.Range("K134").Value = 6988.6158
$ws.Range("L134").Value = 18689.517
$ws.Range("M134").Value = -4453.6158
$ws.Range("N134").Value = -23759.517

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 885.6087
$ws.Range("I107").Value = 315.82352
$ws.Range("K107").Value = 315.82352
$ws.Range("M107").Value = 1604.17648

# Row 122
$ws.Range("H122").Value = 2489.25
$ws.Range("I122").Value = 975
$ws.Range("J122").Value = 4003.5
$ws.Range("K122").Value = 2925
$ws.Range("L122").Value = 12010.5
$ws.Range("M122").Value = -475
$ws.Range("N122").Value = -16910.5

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 132
$ws.Range("H132").Value = 2444.3262
$ws.Range("I132").Value = 1976.0435
$ws.Range("J132").Value = 2912.6086
$ws.Range("K132").Value = 5928.1305
$ws.Range("L132").Value = 8737.825800000001
$ws.Range("M132").Value = -3398.1305
$ws.Range("N132").Value = -13797.8258

# Row 134
$ws.Range("H134").Value = 1587.6389
$ws.Range("I134").Value = 1127
$ws.Range("J134").Value = 2634.5454
$ws.Range("K134").Value = 3381
$ws.Range("L134").Value = 7903.6362
$ws.Range("M134").Value = -846
$ws.Range("N134").Value = -12973.6362

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1061.7021
$ws.Range("I131").Value = 590
$ws.Range("J131").Value = 1144.25
$ws.Range("K131").Value = 1770
$ws.Range("L131").Value = 3432.75
$ws.Range("M131").Value = 3270
$ws.Range("N131").Value = -13512.75

# Row 133
$ws.Range("H133").Value = 3603.9666
$ws.Range("I133").Value = 1895
$ws.Range("J133").Value = 4743.278
$ws.Range("K133").Value = 5685
$ws.Range("L133").Value = 14229.834
$ws.Range("M133").Value = -625
$ws.Range("N133").Value = -24349.834

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("K5").Value = 500
$ws.Range("M5").Value = -388

# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# Row 122
$ws.Range("H122").Value = 3344.2632
$ws.Range("I122").Value = 3033.8125
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9101.4375
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6651.4375
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 3004797.5
$ws.Range("J2").Value = 3004797.5
$ws.Range("L2").Value = 3004797.5
$ws.Range("N2").Value = -3005021.5

# Row 122
$ws.Range("H122").Value = 9600.134
$ws.Range("I122").Value = 10077.077
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 30231.231
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -27781.231
$ws.Range("N122").Value = -24400

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2978.1052
$ws.Range("I122").Value = 3250.6667
$ws.Range("J122").Value = 2852.3076
$ws.Range("K122").Value = 9752.000100000001
$ws.Range("L122").Value = 8556.9228
$ws.Range("M122").Value = -7302.000100000001
$ws.Range("N122").Value = -13456.9228

# Row 123
$ws.Range("H123").Value = 37666.168
$ws.Range("J123").Value = 37666.168
$ws.Range("L123").Value = 37666.168
$ws.Range("N123").Value = -47466.168

# Row 124
$ws.Range("H124").Value = 33113.332
$ws.Range("J124").Value = 33113.332
$ws.Range("L124").Value = 33113.332
$ws.Range("N124").Value = -42933.332

# Row 125
$ws.Range("H125").Value = 29400
$ws.Range("J125").Value = 29400
$ws.Range("L125").Value = 29400
$ws.Range("N125").Value = -39240

# Row 126
$ws.Range("H126").Value = 6070.7
$ws.Range("I126").Value = 6088.375
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 18265.125
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -15795.125
$ws.Range("N126").Value = -22940

# Row 129
$ws.Range("H129").Value = 29730
$ws.Range("J129").Value = 29730
$ws.Range("L129").Value = 29730
$ws.Range("N129").Value = -39730

# Row 130
$ws.Range("H130").Value = 30341.666
$ws.Range("J130").Value = 30341.666
$ws.Range("L130").Value = 30341.666
$ws.Range("N130").Value = -40381.666

# Row 131
$ws.Range("H131").Value = 37907.5
$ws.Range("J131").Value = 37907.5
$ws.Range("L131").Value = 37907.5
$ws.Range("N131").Value = -47987.5

# Row 136
$ws.Range("H136").Value = 589512.5600000001
$ws.Range("I136").Value = 842847.7
$ws.Range("J136").Value = 2841.842
$ws.Range("K136").Value = 2528543.1
$ws.Range("L136").Value = 8525.526
$ws.Range("N136").Value = -13625.526

Write-Output "Edits applied."